# The sheet had a handful of stray empty "placeholder" cells (inline
# strings with no text) scattered through the option columns, plus an
# entire extra question row (row 11) that needs to go away, and the
# sheet view no longer needs to be displayed right-to-left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the scattered empty placeholder cells (J2, I3:I5, G6, E7, I8, G9, G10)
$ws.Range("J2").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("G10").ClearContents()

# Drop the last question (row 11) entirely, shifting nothing else up
# since it is the final row in the used range.
$ws.Rows("11:11").Delete()

# Turn off right-to-left display for the sheet view.
$ws.Application.ActiveWindow.DisplayRightToLeft = $false
